# issue #5: add legislator_id, name, date into dataframe
#
# The "股票" (Stocks) sheet gains three new trailing columns:
#   H = date             ("2012-04-30" for every data row)
#   I = legislator_name  ("楊麗環" for every data row)
#   J = legislator_id    (960 for every data row)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$legislatorDate = "2012-04-30"
$legislatorName = "楊麗環"
$legislatorId = 960

$lastRow = 39

# Header row (row 1): bold/centered header style, matching the existing
# B1:G1 headers. Copy the header format onto the new header cells so they
# pick up the same style as the rest of row 1.
$ws.Cells.Item(1, 8).Value = "date"
$ws.Cells.Item(1, 9).Value = "legislator_name"
$ws.Cells.Item(1, 10).Value = "legislator_id"

$ws.Cells.Item(1, 2).Copy() | Out-Null
$ws.Range("H1:J1").PasteSpecial(-4122) | Out-Null

# Data rows (2..39): plain cell style, matching B2:G39.
# Force the date column to be stored as literal text (not an auto-converted
# Excel date serial) by pre-formatting it as Text before assigning it.
$ws.Range("H2:H" + $lastRow).NumberFormat = "@"

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = $legislatorDate
    $ws.Cells.Item($r, 9).Value = $legislatorName
    $ws.Cells.Item($r, 10).Value = $legislatorId
}

$excel.CutCopyMode = $false
